$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Add the new "Romania" data row (row 31) under the existing urbanization series.
$ws.Range("A31").Value = "Romania"
$ws.Range("B31").Value = 53.9
$ws.Range("C31").Value = 53.887
$ws.Range("D31").Value = 53.9
$ws.Range("E31").Value = 53.936
$ws.Range("F31").Value = 53.998
$ws.Range("G31").Value = 54.084
$ws.Range("H31").Value = 54.194
$ws.Range("I31").Value = 54.329
$ws.Range("J31").Value = 54.489

# Match the author's final selection/scroll position on the Data sheet.
$ws.Activate()
$ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 22
